$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("AI4").Value = 0.458
$ws.Range("AJ4").Value = 0.13
$ws.Range("AK4").Value = 0.361
$ws.Range("AU4").Value = 0.307
$ws.Range("AV4").Value = 0.032
$ws.Range("AW4").Value = 0.178
$ws.Range("BA4").Value = 1.375
$ws.Range("BB4").Value = 0.005
$ws.Range("BC4").Value = 0.072
$ws.Range("BG4").Value = 0.5
$ws.Range("BH4").Value = 0.25
$ws.Range("BI4").Value = 0.5
$ws.Range("BM4").Value = 0.417
$ws.Range("BN4").Value = 0.007
$ws.Range("BO4").Value = 0.083
$ws.Range("BP4").Value = 0.458
$ws.Range("BQ4").Value = 0.765
$ws.Range("E4").Value = 0.583
$ws.Range("F4").Value = 0.068
$ws.Range("G4").Value = 0.26
$ws.Range("N4").Value = 0.542
$ws.Range("O4").Value = 0.075
$ws.Range("P4").Value = 0.273
$ws.Range("Q4").Value = 0.333
$ws.Range("R4").Value = 0.167
$ws.Range("S4").Value = 0.408
$ws.Range("W4").Value = 0.417
$ws.Range("AI5").Value = 0.433
$ws.Range("AJ5").Value = 0.13
$ws.Range("AK5").Value = 0.361
$ws.Range("AU5").Value = 0.567
$ws.Range("AV5").Value = 0.13
$ws.Range("AW5").Value = 0.361
$ws.Range("BA5").Value = 1.133
$ws.Range("BB5").Value = 0.005
$ws.Range("BC5").Value = 0.071
$ws.Range("BG5").Value = 0.317
$ws.Range("BH5").Value = 0.101
$ws.Range("BI5").Value = 0.318
$ws.Range("BM5").Value = 0.383
$ws.Range("BN5").Value = 0.016
$ws.Range("BO5").Value = 0.126
$ws.Range("BP5").Value = 0.378
$ws.Range("BQ5").Value = 0.427
$ws.Range("E5").Value = 0.725
$ws.Range("F5").Value = 0.077
$ws.Range("G5").Value = 0.277
$ws.Range("N5").Value = 0.692
$ws.Range("O5").Value = 0.035
$ws.Range("P5").Value = 0.188
$ws.Range("Q5").Value = 0.133
$ws.Range("R5").Value = 0.02
$ws.Range("S5").Value = 0.141
$ws.Range("W5").Value = 0.383
$ws.Range("X5").Value = 0.141
$ws.Range("Y5").Value = 0.375
$ws.Range("AI6").Value = 0.445
$ws.Range("AU6").Value = 0.398
$ws.Range("BA6").Value = 1.232
$ws.Range("BG6").Value = 0.388
$ws.Range("BM6").Value = 0.399
$ws.Range("BP6").Value = 0.411
$ws.Range("BQ6").Value = 0.545
$ws.Range("E6").Value = 0.646
$ws.Range("N6").Value = 0.608
$ws.Range("Q6").Value = 0.19
$ws.Range("W6").Value = 0.399
$ws.Range("AI7").Value = 0.438
$ws.Range("AU7").Value = 0.485
$ws.Range("BA7").Value = 1.169
$ws.Range("BG7").Value = 0.342
$ws.Range("BM7").Value = 0.389
$ws.Range("BP7").Value = 0.39
$ws.Range("BQ7").Value = 0.467
$ws.Range("E7").Value = 0.691
$ws.Range("N7").Value = 0.656
$ws.Range("Q7").Value = 0.151
$ws.Range("W7").Value = 0.389
$ws.Range("AI8").Value = 0.621
$ws.Range("AJ8").Value = 0.136
$ws.Range("AK8").Value = 0.368
$ws.Range("AU8").Value = 0.511
$ws.Range("AV8").Value = 0.107
$ws.Range("AW8").Value = 0.326
$ws.Range("BA8").Value = 1.502
$ws.Range("BB8").Value = 0.046
$ws.Range("BC8").Value = 0.215
$ws.Range("BG8").Value = 0.438
$ws.Range("BH8").Value = 0.194
$ws.Range("BI8").Value = 0.44
$ws.Range("BM8").Value = 0.472
$ws.Range("BN8").Value = 0.085
$ws.Range("BO8").Value = 0.292
$ws.Range("BP8").Value = 0.501
$ws.Range("BQ8").Value = 0.6
$ws.Range("E8").Value = 0.716
$ws.Range("F8").Value = 0.118
$ws.Range("G8").Value = 0.344
$ws.Range("N8").Value = 0.833
$ws.Range("O8").Value = 0.033
$ws.Range("P8").Value = 0.18
$ws.Range("Q8").Value = 0.196
$ws.Range("R8").Value = 0.08
$ws.Range("S8").Value = 0.282
$ws.Range("W8").Value = 0.555
$ws.Range("X8").Value = 0.105
$ws.Range("Y8").Value = 0.323
$ws.Range("AI9").Value = 0.75
$ws.Range("AJ9").Value = 0.188
$ws.Range("AK9").Value = 0.433
$ws.Range("BA9").Value = 1.5
$ws.Range("BM9").Value = 0.5
$ws.Range("BN9").Value = 0.25
$ws.Range("BO9").Value = 0.5
$ws.Range("BP9").Value = 0.5
$ws.Range("BQ9").Value = 0.622
$ws.Range("E9").Value = 0.75
$ws.Range("F9").Value = 0.188
$ws.Range("G9").Value = 0.433
$ws.Range("N9").Value = 0.75
$ws.Range("O9").Value = 0.188
$ws.Range("P9").Value = 0.433
$ws.Range("AI10").Value = 0.75
$ws.Range("AJ10").Value = 0.188
$ws.Range("AK10").Value = 0.433
$ws.Range("BA10").Value = 1.75
$ws.Range("BB10").Value = 0.188
$ws.Range("BC10").Value = 0.433
$ws.Range("BM10").Value = 0.5
$ws.Range("BN10").Value = 0.25
$ws.Range("BO10").Value = 0.5
$ws.Range("BP10").Value = 0.583
$ws.Range("BQ10").Value = 0.711
$ws.Range("E10").Value = 0.75
$ws.Range("F10").Value = 0.188
$ws.Range("G10").Value = 0.433
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("W10").Value = 0.75
$ws.Range("X10").Value = 0.188
$ws.Range("Y10").Value = 0.433
$ws.Range("AI11").Value = 0.75
$ws.Range("AJ11").Value = 0.188
$ws.Range("AK11").Value = 0.433
$ws.Range("AU11").Value = 0.5
$ws.Range("AV11").Value = 0.25
$ws.Range("AW11").Value = 0.5
$ws.Range("BA11").Value = 1.75
$ws.Range("BB11").Value = 0.188
$ws.Range("BC11").Value = 0.433
$ws.Range("BM11").Value = 0.5
$ws.Range("BN11").Value = 0.25
$ws.Range("BO11").Value = 0.5
$ws.Range("BP11").Value = 0.583
$ws.Range("BQ11").Value = 0.711
$ws.Range("E11").Value = 0.75
$ws.Range("F11").Value = 0.188
$ws.Range("G11").Value = 0.433
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("W11").Value = 0.75
$ws.Range("X11").Value = 0.188
$ws.Range("Y11").Value = 0.433
$ws.Range("AV12").Value = 6
$ws.Range("AW12").Value = 2.449
$ws.Range("BA12").Value = 3.333
$ws.Range("BB12").Value = 0.222
$ws.Range("BC12").Value = 0.471
$ws.Range("BP12").Value = 1.111
$ws.Range("BQ12").Value = 1.139
$ws.Range("N12").Value = 1.25
$ws.Range("O12").Value = 0.188
$ws.Range("P12").Value = 0.433
$ws.Range("W12").Value = 1.333
$ws.Range("X12").Value = 0.222
$ws.Range("Y12").Value = 0.471
$ws.Range("BP13").Value = 0.834
$ws.Range("BQ13").Value = 0.6

Write-Host "Applied 182 cell changes"
